# Apply cell value updates per the diff (odds data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.75  # G2: 1.8 -> 1.75
$ws.Cells.Item(2, 9).Value = 5.25  # I2: 5 -> 5.25
$ws.Cells.Item(2, 12).Value = 6  # L2: 5.5 -> 6
$ws.Cells.Item(2, 13).Value = 1.11  # M2: 1.1 -> 1.11
$ws.Cells.Item(2, 14).Value = 6.5  # N2: 7 -> 6.5
$ws.Cells.Item(2, 17).Value = 1.93  # Q2: 1.9 -> 1.93
$ws.Cells.Item(2, 18).Value = 1.93  # R2: 1.95 -> 1.93
$ws.Cells.Item(2, 25).Value = 1.57  # Y2: 1.53 -> 1.57
$ws.Cells.Item(2, 26).Value = 2.25  # Z2: 2.38 -> 2.25
$ws.Cells.Item(2, 29).Value = 5  # AC2: 5.5 -> 5
$ws.Cells.Item(2, 44).Value = 51  # AR2: 41 -> 51
# Row 3
$ws.Cells.Item(3, 7).Value = 2.4  # G3: 2.5 -> 2.4
$ws.Cells.Item(3, 9).Value = 3  # I3: 2.9 -> 3
$ws.Cells.Item(3, 10).Value = 3.25  # J3: 3.4 -> 3.25
$ws.Cells.Item(3, 11).Value = 1.91  # K3: 1.95 -> 1.91
$ws.Cells.Item(3, 12).Value = 4  # L3: 3.75 -> 4
$ws.Cells.Item(3, 17).Value = 1.93  # Q3: 1.9 -> 1.93
$ws.Cells.Item(3, 18).Value = 1.93  # R3: 1.95 -> 1.93
$ws.Cells.Item(3, 21).Value = 3.95  # U3: 3.9 -> 3.95
$ws.Cells.Item(3, 22).Value = 1.24  # V3: 1.25 -> 1.24
$ws.Cells.Item(3, 30).Value = 10  # AD3: 11 -> 10
$ws.Cells.Item(3, 31).Value = 10  # AE3: 11 -> 10
$ws.Cells.Item(3, 40).Value = 7.5  # AN3: 7 -> 7.5
$ws.Cells.Item(3, 43).Value = 34  # AQ3: 29 -> 34
# Row 4
$ws.Cells.Item(4, 7).Value = 2.8  # G4: 2.75 -> 2.8
$ws.Cells.Item(4, 9).Value = 2.63  # I4: 2.7 -> 2.63
$ws.Cells.Item(4, 10).Value = 3.75  # J4: 3.6 -> 3.75
$ws.Cells.Item(4, 12).Value = 3.5  # L4: 3.6 -> 3.5
$ws.Cells.Item(4, 13).Value = 1.11  # M4: 1.1 -> 1.11
$ws.Cells.Item(4, 14).Value = 6.5  # N4: 7 -> 6.5
$ws.Cells.Item(4, 31).Value = 12  # AE4: 11 -> 12
$ws.Cells.Item(4, 33).Value = 29  # AG4: 26 -> 29
$ws.Cells.Item(4, 43).Value = 26  # AQ4: 29 -> 26
# Row 5
$ws.Cells.Item(5, 7).Value = 2.8  # G5: 2.92 -> 2.8
$ws.Cells.Item(5, 8).Value = 2.5  # H5: 2.47 -> 2.5
$ws.Cells.Item(5, 9).Value = 3.1  # I5: 3 -> 3.1
$ws.Cells.Item(5, 10).Value = 3.6  # J5: 3.75 -> 3.6
$ws.Cells.Item(5, 11).Value = 1.75  # K5: 1.72 -> 1.75
$ws.Cells.Item(5, 12).Value = 3.9  # L5: 3.85 -> 3.9
$ws.Cells.Item(5, 14).Value = 4.45  # N5: 4.4 -> 4.45
$ws.Cells.Item(5, 19).Value = 2.92  # S5: 2.95 -> 2.92
$ws.Cells.Item(5, 20).Value = 1.36  # T5: 1.35 -> 1.36
$ws.Cells.Item(5, 23).Value = 5.3  # W5: 5.4 -> 5.3
$ws.Cells.Item(5, 26).Value = 2.05  # Z5: 2.02 -> 2.05
$ws.Cells.Item(5, 29).Value = 6  # AC5: 6.1 -> 6
$ws.Cells.Item(5, 30).Value = 12.5  # AD5: 13.5 -> 12.5
$ws.Cells.Item(5, 31).Value = 11  # AE5: 11.25 -> 11
$ws.Cells.Item(5, 32).Value = 37  # AF5: 40 -> 37
$ws.Cells.Item(5, 33).Value = 32  # AG5: 35 -> 32
$ws.Cells.Item(5, 35).Value = 4.45  # AI5: 4.4 -> 4.45
$ws.Cells.Item(5, 36).Value = 5.2  # AJ5: 5.1 -> 5.2
$ws.Cells.Item(5, 40).Value = 6.3  # AN5: 6.1 -> 6.3
$ws.Cells.Item(5, 41).Value = 14  # AO5: 13.5 -> 14
$ws.Cells.Item(5, 42).Value = 12  # AP5: 11.5 -> 12
$ws.Cells.Item(5, 44).Value = 40  # AR5: 37 -> 40
# Row 6
$ws.Cells.Item(6, 7).Value = 3.6  # G6: 3.8 -> 3.6
$ws.Cells.Item(6, 8).Value = 2.85  # H6: 2.87 -> 2.85
$ws.Cells.Item(6, 9).Value = 2.22  # I6: 2.15 -> 2.22
$ws.Cells.Item(6, 10).Value = 4.2  # J6: 4.35 -> 4.2
$ws.Cells.Item(6, 12).Value = 2.87  # L6: 2.8 -> 2.87
$ws.Cells.Item(6, 23).Value = 4.3  # W6: 4.35 -> 4.3
$ws.Cells.Item(6, 27).Value = 2  # AA6: 2.05 -> 2
$ws.Cells.Item(6, 28).Value = 1.72  # AB6: 1.7 -> 1.72
$ws.Cells.Item(6, 30).Value = 18  # AD6: 19 -> 18
$ws.Cells.Item(6, 31).Value = 12.5  # AE6: 13 -> 12.5
$ws.Cells.Item(6, 32).Value = 55  # AF6: 60 -> 55
$ws.Cells.Item(6, 33).Value = 40  # AG6: 45 -> 40
$ws.Cells.Item(6, 34).Value = 50  # AH6: 55 -> 50
$ws.Cells.Item(6, 36).Value = 5.6  # AJ6: 5.7 -> 5.6
$ws.Cells.Item(6, 37).Value = 16.5  # AK6: 17 -> 16.5
$ws.Cells.Item(6, 38).Value = 100  # AL6: 110 -> 100
$ws.Cells.Item(6, 40).Value = 5.9  # AN6: 5.7 -> 5.9
$ws.Cells.Item(6, 41).Value = 9.5  # AO6: 9 -> 9.5
$ws.Cells.Item(6, 42).Value = 9.25  # AP6: 9 -> 9.25
$ws.Cells.Item(6, 43).Value = 22  # AQ6: 21 -> 22
$ws.Cells.Item(6, 45).Value = 37  # AS6: 40 -> 37
# Row 7
$ws.Cells.Item(7, 9).Value = 3.3  # I7: 3.2 -> 3.3
$ws.Cells.Item(7, 13).Value = 1.13  # M7: 1.11 -> 1.13
$ws.Cells.Item(7, 14).Value = 6  # N7: 6.5 -> 6
$ws.Cells.Item(7, 15).Value = 1.57  # O7: 1.53 -> 1.57
$ws.Cells.Item(7, 16).Value = 2.38  # P7: 2.5 -> 2.38
$ws.Cells.Item(7, 17).Value = 2.1  # Q7: 2.03 -> 2.1
$ws.Cells.Item(7, 18).Value = 1.78  # R7: 1.83 -> 1.78
$ws.Cells.Item(7, 40).Value = 7.5  # AN7: 7 -> 7.5
# Row 8
$ws.Cells.Item(8, 19).Value = 2.08  # S8: 2.07 -> 2.08
$ws.Cells.Item(8, 20).Value = 1.73  # T8: 1.69 -> 1.73
# Row 9
$ws.Cells.Item(9, 20).Value = 1.67  # T9: 1.63 -> 1.67
# Row 10
$ws.Cells.Item(10, 12).Value = 2.63  # L10: 2.62 -> 2.63
$ws.Cells.Item(10, 20).Value = 1.67  # T10: 1.63 -> 1.67
# Row 11
$ws.Cells.Item(11, 20).Value = 1.7  # T11: 1.67 -> 1.7
# Row 12
$ws.Cells.Item(12, 10).Value = 2.88  # J12: 2.87 -> 2.88
$ws.Cells.Item(12, 20).Value = 1.75  # T12: 1.72 -> 1.75
# Row 13
$ws.Cells.Item(13, 7).Value = 1.38  # G13: 1.4 -> 1.38
$ws.Cells.Item(13, 8).Value = 4.5  # H13: 4.33 -> 4.5
$ws.Cells.Item(13, 9).Value = 8  # I13: 7.5 -> 8
$ws.Cells.Item(13, 10).Value = 1.91  # J13: 1.87 -> 1.91
$ws.Cells.Item(13, 14).Value = 12  # N13: 13 -> 12
$ws.Cells.Item(13, 19).Value = 1.7  # S13: 1.67 -> 1.7
$ws.Cells.Item(13, 30).Value = 6.5  # AD13: 7 -> 6.5
$ws.Cells.Item(13, 40).Value = 21  # AN13: 19 -> 21
